$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A2's format (bold font, border, alignment) down to A3 and A4
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 1
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 1
